$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record was inserted at row 28 (shifting the existing rows
# 28..95 down to 29..96). The new row re-uses the same market / product
# metadata as the old row 28, but with a new date and volume.
$ws.Rows.Item(28).Insert()

$ws.Cells.Item(28,1).Value = 5
$ws.Cells.Item(28,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(28,3).Value = "Maule"
$ws.Cells.Item(28,4).Value = "2023-06-05"
$ws.Cells.Item(28,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(28,5).Value = 7
$ws.Cells.Item(28,6).Value = 100112040
$ws.Cells.Item(28,7).Value = "Cilantro"
$ws.Cells.Item(28,8).Value = "Sin especificar"
$ws.Cells.Item(28,9).Value = "Primera"
$ws.Cells.Item(28,10).Value = 300
$ws.Cells.Item(28,11).Value = 7000
$ws.Cells.Item(28,12).Value = 7000
$ws.Cells.Item(28,13).Value = 7000
$ws.Cells.Item(28,14).Value = "`$/caja 36 atados"
$ws.Cells.Item(28,15).Value = "Región del Maule"
$ws.Cells.Item(28,16).Value = 194
$ws.Cells.Item(28,17).Value = 36
$ws.Cells.Item(28,18).Value = "Hortaliza"
